# Update the "Förändrad" date column (C) for all data rows (2-45) from
# 2023-09-02 (serial 45171) to 2023-09-03 (serial 45172).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
